# BalanceDaily.xlsx - "Se agregan nuevos bots"
# Appends a new daily balance row (2021-02-06) to the BalanceDaily table on Hoja1,
# which grows the table/autofilter range from A1:D42 to A1:D43 and moves the
# active selection to B44 (the next empty cell below the table).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Expand the Excel Table (ListObject) by one row so the table ref / autoFilter
# ref grow from A1:D42 to A1:D43, matching the structured-table behaviour.
$lo = $ws.ListObjects.Item("BalanceDaily")
$null = $lo.ListRows.Add()

# Fill in the new row's data. Column C/D keep the table's calculated-column
# formulas (same pattern as every other data row).
$ws.Range("A43").Value = 44233
$ws.Range("B43").Value = 0.0086053699999999993
$ws.Range("C43").Formula = "=ROUND(IFERROR(BalanceDaily[[#This Row],[ValueLTC]]-B42,0),8)"
$ws.Range("D43").Formula = "=BalanceDaily[[#This Row],[IncrementDaily]]/24"

# Move the selection to the next empty cell under the table, like Excel does
# after you finish typing a new row.
$null = $ws.Range("B44").Select()
